# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.710.90'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.628.86'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.08'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.254'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.44'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.854.52'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').Value = '1.627.77'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.554'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').Value = '0.0₃0759'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.99'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').Value = '25.727.09'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.42'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.14'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.87'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.23'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.81'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.55'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.43'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.38'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '1.130.66'
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.540'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.04'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.50'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Value = '1.763.98'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0511'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('E51').Value = '  +3.20%  '
